$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.290.74'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '2.602.07'
$ws.Range('E3').Value = '  +2.21%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''306.88'
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = '''99.40'
$ws.Range('E6').Value = '  -3.63%  '
$ws.Range('E7').Value = '  -1.32%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '''0.577'
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').Value = '''39.30'
$ws.Range('E10').Value = '  +0.91%  '
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').Value = '''54.09'
$ws.Range('E11').Value = '  -1.17%  '
$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D12').Value = '''0.0840'
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('D13').Value = '''8.07'
$ws.Range('E13').Value = '  +1.21%  '
$ws.Range('D14').Value = '3.001.71'
$ws.Range('E14').Value = '  +2.37%  '
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('D16').Value = '2.610.37'
$ws.Range('E16').Value = '  +2.81%  '
$ws.Range('D17').Value = '''0.914'
$ws.Range('E17').Value = '  +1.80%  '
$ws.Range('D18').Value = '''14.87'
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('D19').Value = '46.372.76'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('E20').Value = '  +0.95%  '
$ws.Range('E21').Value = '  -7.12%  '
$ws.Range('D22').Value = '''6.68'
$ws.Range('E22').Value = '  +0.70%  '
$ws.Range('D23').Value = '''71.16'
$ws.Range('E23').Value = '  +1.64%  '
$ws.Range('D24').Value = '''271.54'
$ws.Range('E24').Value = '  +6.60%  '
$ws.Range('D25').Value = '''3.02'
$ws.Range('E25').Value = '  +1.09%  '
$ws.Range('D26').Value = '''2.15'
$ws.Range('E26').Value = '  +1.09%  '
$ws.Range('D27').Value = '''29.19'
$ws.Range('E27').Value = '  +20.91%  '
$ws.Range('E29').Value = '  -0.80%  '
$ws.Range('D30').Value = '''10.54'
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = '''2.27'
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D32').Value = '''38.39'
$ws.Range('E32').Value = '  -9.18%  '
$ws.Range('D33').Value = '''6.30'
$ws.Range('E33').Value = '  +4.82%  '
$ws.Range('E34').Value = '  -4.47%  '
$ws.Range('E35').Value = '  -2.26%  '
$ws.Range('D36').Value = '''2.22'
$ws.Range('E36').Value = '  +1.87%  '
$ws.Range('D37').Value = '''0.0833'
$ws.Range('E37').Value = '  -1.73%  '
$ws.Range('D38').Value = '''151.22'
$ws.Range('E38').Value = '  +0.63%  '
$ws.Range('E39').Value = '  +3.21%  '
$ws.Range('E40').Value = '  +1.02%  '
$ws.Range('D41').Value = '''23.05'
$ws.Range('E41').Value = '  +28.56%  '
$ws.Range('D42').Value = '''15.77'
$ws.Range('E42').Value = '  -5.78%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '''0.0328'
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').Value = '''3.59'
$ws.Range('E44').Value = '  +0.70%  '
$ws.Range('D45').Value = '''4.04'
$ws.Range('E45').Value = '  -4.54%  '
$ws.Range('D46').Value = '2.113.81'
$ws.Range('E46').Value = '  +5.53%  '
$ws.Range('E47').Value = '  -0.18%  '
$ws.Range('D48').Value = '''93.01'
$ws.Range('E48').Value = '  -1.01%  '
$ws.Range('D49').Value = '''9.54'
$ws.Range('E49').Value = '  +6.93%  '
$ws.Range('E50').Value = '  -5.20%  '
$ws.Range('D51').Value = '''108.18'
$ws.Range('E51').Value = '  +0.67%  '
